# SSDM-13256: Fixing failing tests.
# Add a "Validation script" column (F) to the sample type export sheet,
# with an example value of "test.py" for the ENTRY type, and move the
# active selection to the new header cell (F2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell (row 2, bold style matches the rest of the header row)
$ws.Range("F2").Value = "Validation script"
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)

# Data cell (row 3, regular style matches the rest of the ENTRY row)
$ws.Range("F3").Value = "test.py"
$ws.Range("E3").Copy()
$ws.Range("F3").PasteSpecial(-4122)

# Drop the now-stray empty styled cells in column J for rows 2 and 3
$ws.Cells.Item(2, 10).ClearContents()
$ws.Cells.Item(3, 10).ClearContents()

# Move the active cell / selection to the new header cell
$ws.Range("F2").Select()
